$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New primer data (replaces rows 1-32, gene renamed bdhA_del/acoA_del/pta_del/ldh_del -> BSU_000NN)
$data = @(
    @("BSU_00020-1", "ACGTCTCACATATGACTAGTGCATGCGATCGGCCATTATG"),
    @("BSU_00020-2", "AGTTCGAACTTAAGAGATCTTTTTATCCTCCTAACGGATAATGT"),
    @("BSU_00020-3", "AGGGTAATGCATGCCTGCAGATTAATACATTATCCGTTAGGAGGATAAAATCCGATACACTGCTGCCGA"),
    @("BSU_00020-4", "GGGGATCCTCTAGAGTCGACCGAGCTGTGCTGTAAACTGC"),
    @("BSU_00020-7", "AGATCTCTTAAGTTCGAACTCGAGC"),
    @("BSU_00020-8", "CTGCAGGCATGCATTACCCT"),
    @("test-BSU_00020-1", "TGTCCACATGTGGATAGGCTG"),
    @("test-BSU_00020-2", "AGCTACACGCTGTCTTGCTTC"),
    @("BSU_00030-1", "ACGTCTCACATATGACTAGTGCATCGTTTTACAGGCTCGC"),
    @("BSU_00030-2", "AGTTCGAACTTAAGAGATCTTATATCGACCTCTTTCAAATATCAC"),
    @("BSU_00030-3", "AGGGTAATGCATGCCTGCAGTATTAGTGATATTTGAAAGAGGTCGATATAAGCGGGTGACACTGATTGTA"),
    @("BSU_00030-4", "GGGGATCCTCTAGAGTCGACCACGGCCTTGGATCGTATGA"),
    @("BSU_00030-7", "AGATCTCTTAAGTTCGAACTCGAGC"),
    @("BSU_00030-8", "CTGCAGGCATGCATTACCCT"),
    @("test-BSU_00030-1", "TCCGATACACTGCTGCCGA"),
    @("test-BSU_00030-2", "AGCTACACGCTGTCTTGCTTC"),
    @("BSU_00040-1", "ACGTCTCACATATGACTAGTCACCTCAGAAACACGCCCTA"),
    @("BSU_00040-2", "AGTTCGAACTTAAGAGATCTTCAGTGTCACCCGCTTTAATTG"),
    @("BSU_00040-3", "AGGGTAATGCATGCCTGCAGCAAGTCGTCAATTAAAGCGGGTGACACTGAAGAAATGAGGTGAGCAATTGT"),
    @("BSU_00040-4", "GGGGATCCTCTAGAGTCGACCGGCTCTTCATGGACAACCT"),
    @("BSU_00040-7", "AGATCTCTTAAGTTCGAACTCGAGC"),
    @("BSU_00040-8", "CTGCAGGCATGCATTACCCT"),
    @("test-BSU_00040-1", "AAGCGAGCATGAAGTGCTTG"),
    @("test-BSU_00040-2", "AGCTACACGCTGTCTTGCTTC"),
    @("BSU_00050-1", "ACGTCTCACATATGACTAGTTCGCACCGGACATCAAATGA"),
    @("BSU_00050-2", "AGTTCGAACTTAAGAGATCTTTGCTCACCTCATTTCTTCACT"),
    @("BSU_00050-3", "AGGGTAATGCATGCCTGCAGGTTAGTGAAGTGAAGAAATGAGGTGAGCAAAAATTTTTTATCACGAATATATCGT"),
    @("BSU_00050-4", "GGGGATCCTCTAGAGTCGACCTGTCAGCCCTTCCCTTACG"),
    @("BSU_00050-7", "AGATCTCTTAAGTTCGAACTCGAGC"),
    @("BSU_00050-8", "CTGCAGGCATGCATTACCCT"),
    @("test-BSU_00050-1", "TCCAAGGCCGTGTACAAACG"),
    @("test-BSU_00050-2", "AGCTACACGCTGTCTTGCTTC")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}

# Remove the old rows 33-56 (alsSD_P43_sub / amyE::P43_budC_ins blocks) entirely
$lastRow = $data.Length + 1
$ws.Range("A" + $lastRow + ":B56").EntireRow.Delete() | Out-Null
